# Generated PowerShell COM-interop script
# Implements: insert two new localization entries ("0e1265aa..." and
# "2f66370b...") ahead of the existing "cc3cdaa2..." row on every sheet,
# marking them "Ready for handoff", and refresh all hyperlinks.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1 (sheet1) ----
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows before the existing "cc3cdaa2..." row (row 5),
# pushing it down to row 7 and opening rows 5-6 for the new entries.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Row 5
$ws.Range("A5").Value = "0e1265aa-4f8b-44ef-8e83-66438e852d0d.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "2016-03-24 16:47:47"

# Row 6
$ws.Range("A6").Value = "2f66370b-b656-486f-b385-4db18e61c7c7.md"
$ws.Range("B6").Value = "Ready for handoff"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "2016-03-24 16:47:47"

# Row 7
$ws.Range("A7").Value = "cc3cdaa2-a023-4e07-9164-a680c6834b7d.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "2016-03-24 16:39:55"

# Hyperlinks are not shifted by row-insert in this engine, so clear them
# all and re-add in the correct final row order.
$ws.Range("A1").Hyperlinks.Delete()

$hyperlinks_sheet1 = @(
    ,@("A2", "3edcf9a3-0613-4e5b-844c-20befde9455d.md")
    ,@("A3", "8aba6dc8-6b45-434b-b38c-f2765e82d759.md")
    ,@("A4", "d66e2f39-0931-4909-b497-ea4e0675319e.md")
    ,@("A5", "0e1265aa-4f8b-44ef-8e83-66438e852d0d.md")
    ,@("A6", "2f66370b-b656-486f-b385-4db18e61c7c7.md")
    ,@("A7", "cc3cdaa2-a023-4e07-9164-a680c6834b7d.md")
)
foreach ($item in $hyperlinks_sheet1) {
    $cellRef = $item[0]
    $displayText = $item[1]
    $address = "https://example.com/" + $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText) | Out-Null
}

# ---- Sheet 2 (sheet2) ----
$ws = $wb.Worksheets.Item(2)

# Insert two blank rows before the existing "cc3cdaa2..." row (row 5),
# pushing it down to row 7 and opening rows 5-6 for the new entries.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Row 5
$ws.Range("A5").Value = "0e1265aa-4f8b-44ef-8e83-66438e852d0d.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "0e1265aa-4f8b-44ef-8e83-66438e852d0d.a23302d567891f7afdbaace2e00cb354916ba03d.zh-cn.xlf"
$ws.Range("E5").Value = "2016-03-24 16:47:42"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("J5").Value = "Include"

# Row 6
$ws.Range("A6").Value = "2f66370b-b656-486f-b385-4db18e61c7c7.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "2f66370b-b656-486f-b385-4db18e61c7c7.cc6f9dc275b784fe7608e5b8012a5b2c8a2eb78c.zh-cn.xlf"
$ws.Range("E6").Value = "2016-03-24 16:47:42"
$ws.Range("H6").Value = "0001-01-01 00:00:00"
$ws.Range("J6").Value = "Include"

# Row 7
$ws.Range("A7").Value = "cc3cdaa2-a023-4e07-9164-a680c6834b7d.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "cc3cdaa2-a023-4e07-9164-a680c6834b7d.6518ae3aa56c44c1d0736a81625548e39f1f49de.zh-cn.xlf"
$ws.Range("E7").Value = "2016-03-24 16:39:51"
$ws.Range("H7").Value = "0001-01-01 00:00:00"
$ws.Range("J7").Value = "Include"

# Hyperlinks are not shifted by row-insert in this engine, so clear them
# all and re-add in the correct final row order.
$ws.Range("A1").Hyperlinks.Delete()

$hyperlinks_sheet2 = @(
    ,@("A2", "3edcf9a3-0613-4e5b-844c-20befde9455d.md")
    ,@("D2", "3edcf9a3-0613-4e5b-844c-20befde9455d.904a83b41c482afd675a64afd621ecd6b3cc3d89.zh-cn.xlf")
    ,@("F2", "3edcf9a3-0613-4e5b-844c-20befde9455d.md")
    ,@("G2", "3edcf9a3-0613-4e5b-844c-20befde9455d.904a83b41c482afd675a64afd621ecd6b3cc3d89.zh-cn.xlf")
    ,@("A3", "8aba6dc8-6b45-434b-b38c-f2765e82d759.md")
    ,@("D3", "8aba6dc8-6b45-434b-b38c-f2765e82d759.ad3672ac52ff8c1c01ad8c29d8961de8b77a39ff.zh-cn.xlf")
    ,@("A4", "d66e2f39-0931-4909-b497-ea4e0675319e.md")
    ,@("D4", "d66e2f39-0931-4909-b497-ea4e0675319e.39c8048b6e853cdd2e1e1621d7e0211b9027a102.zh-cn.xlf")
    ,@("A5", "0e1265aa-4f8b-44ef-8e83-66438e852d0d.md")
    ,@("D5", "0e1265aa-4f8b-44ef-8e83-66438e852d0d.a23302d567891f7afdbaace2e00cb354916ba03d.zh-cn.xlf")
    ,@("A6", "2f66370b-b656-486f-b385-4db18e61c7c7.md")
    ,@("D6", "2f66370b-b656-486f-b385-4db18e61c7c7.cc6f9dc275b784fe7608e5b8012a5b2c8a2eb78c.zh-cn.xlf")
    ,@("A7", "cc3cdaa2-a023-4e07-9164-a680c6834b7d.md")
    ,@("D7", "cc3cdaa2-a023-4e07-9164-a680c6834b7d.6518ae3aa56c44c1d0736a81625548e39f1f49de.zh-cn.xlf")
)
foreach ($item in $hyperlinks_sheet2) {
    $cellRef = $item[0]
    $displayText = $item[1]
    $address = "https://example.com/" + $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText) | Out-Null
}

# ---- Sheet 3 (sheet3) ----
$ws = $wb.Worksheets.Item(3)

# Insert two blank rows before the existing "cc3cdaa2..." row (row 5),
# pushing it down to row 7 and opening rows 5-6 for the new entries.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# Row 5
$ws.Range("A5").Value = "0e1265aa-4f8b-44ef-8e83-66438e852d0d.md"
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "0e1265aa-4f8b-44ef-8e83-66438e852d0d.a23302d567891f7afdbaace2e00cb354916ba03d.de-de.xlf"
$ws.Range("E5").Value = "2016-03-24 16:47:47"
$ws.Range("H5").Value = "0001-01-01 00:00:00"
$ws.Range("J5").Value = "Include"

# Row 6
$ws.Range("A6").Value = "2f66370b-b656-486f-b385-4db18e61c7c7.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "2f66370b-b656-486f-b385-4db18e61c7c7.cc6f9dc275b784fe7608e5b8012a5b2c8a2eb78c.de-de.xlf"
$ws.Range("E6").Value = "2016-03-24 16:47:47"
$ws.Range("H6").Value = "0001-01-01 00:00:00"
$ws.Range("J6").Value = "Include"

# Row 7
$ws.Range("A7").Value = "cc3cdaa2-a023-4e07-9164-a680c6834b7d.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "cc3cdaa2-a023-4e07-9164-a680c6834b7d.6518ae3aa56c44c1d0736a81625548e39f1f49de.de-de.xlf"
$ws.Range("E7").Value = "2016-03-24 16:39:55"
$ws.Range("H7").Value = "0001-01-01 00:00:00"
$ws.Range("J7").Value = "Include"

# Hyperlinks are not shifted by row-insert in this engine, so clear them
# all and re-add in the correct final row order.
$ws.Range("A1").Hyperlinks.Delete()

$hyperlinks_sheet3 = @(
    ,@("A2", "3edcf9a3-0613-4e5b-844c-20befde9455d.md")
    ,@("D2", "3edcf9a3-0613-4e5b-844c-20befde9455d.904a83b41c482afd675a64afd621ecd6b3cc3d89.de-de.xlf")
    ,@("F2", "3edcf9a3-0613-4e5b-844c-20befde9455d.md")
    ,@("G2", "3edcf9a3-0613-4e5b-844c-20befde9455d.904a83b41c482afd675a64afd621ecd6b3cc3d89.de-de.xlf")
    ,@("A3", "8aba6dc8-6b45-434b-b38c-f2765e82d759.md")
    ,@("D3", "8aba6dc8-6b45-434b-b38c-f2765e82d759.ad3672ac52ff8c1c01ad8c29d8961de8b77a39ff.de-de.xlf")
    ,@("A4", "d66e2f39-0931-4909-b497-ea4e0675319e.md")
    ,@("D4", "d66e2f39-0931-4909-b497-ea4e0675319e.39c8048b6e853cdd2e1e1621d7e0211b9027a102.de-de.xlf")
    ,@("A5", "0e1265aa-4f8b-44ef-8e83-66438e852d0d.md")
    ,@("D5", "0e1265aa-4f8b-44ef-8e83-66438e852d0d.a23302d567891f7afdbaace2e00cb354916ba03d.de-de.xlf")
    ,@("A6", "2f66370b-b656-486f-b385-4db18e61c7c7.md")
    ,@("D6", "2f66370b-b656-486f-b385-4db18e61c7c7.cc6f9dc275b784fe7608e5b8012a5b2c8a2eb78c.de-de.xlf")
    ,@("A7", "cc3cdaa2-a023-4e07-9164-a680c6834b7d.md")
    ,@("D7", "cc3cdaa2-a023-4e07-9164-a680c6834b7d.6518ae3aa56c44c1d0736a81625548e39f1f49de.de-de.xlf")
)
foreach ($item in $hyperlinks_sheet3) {
    $cellRef = $item[0]
    $displayText = $item[1]
    $address = "https://example.com/" + $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $displayText) | Out-Null
}
